$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Sieve Bootstrap"
$ws.Range("B2").Value = "3/10"
$ws.Range("C2").Value = 76.8
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 0.5569108207977537

# Row 3
$ws.Range("A3").Value = "Block Bootstrapping"
$ws.Range("B3").Value = "2/10"
$ws.Range("C3").Value = 51.2
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 4.001127937737089

# Row 4
$ws.Range("A4").Value = "AREPD"
$ws.Range("B4").Value = "1/10"
$ws.Range("C4").Value = 25.6
$ws.Range("D4").Value = 20
$ws.Range("E4").Value = 3.587937739237214

# Rows 5-10 keep their model names and numeric data, only the
# "Comparaciones_Significativas" column text is refreshed to "0/10".
$ws.Range("B5").Value = "0/10"
$ws.Range("B6").Value = "0/10"
$ws.Range("B7").Value = "0/10"
$ws.Range("B8").Value = "0/10"
$ws.Range("B9").Value = "0/10"
$ws.Range("B10").Value = "0/10"

$wb.Save()
